$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1275699870682362
$ws.Range("C2").Value = 0.2949975552939468
$ws.Range("D2").Value = 0.150653331260784
$ws.Range("E2").Value = 0.3881408652290866
$ws.Range("F2").Value = 0.3748164570332126
$ws.Range("G2").Value = 23

$ws.Range("B3").Value = -0.01521116492496113
$ws.Range("C3").Value = 0.2722868036881911
$ws.Range("D3").Value = 0.1278291975161459
$ws.Range("E3").Value = 0.3575320929876729
$ws.Range("F3").Value = 0.3656144205965373
$ws.Range("G3").Value = 22

$ws.Range("B4").Value = 0.1098904018663597
$ws.Range("C4").Value = 0.2824368076085529
$ws.Range("D4").Value = 0.1273013017406104
$ws.Range("E4").Value = 0.3567930797263457
$ws.Range("F4").Value = 0.3478313835526827
$ws.Range("G4").Value = 21

$ws.Range("B5").Value = 0.0325328965656304
$ws.Range("C5").Value = 0.2428801592263986
$ws.Range("D5").Value = 0.08899619181356792
$ws.Range("E5").Value = 0.2983222952002882
$ws.Range("F5").Value = 0.3042467877348358
$ws.Range("G5").Value = 20

$ws.Range("B6").Value = 0.07522128658533693
$ws.Range("C6").Value = 0.2678359888536521
$ws.Range("D6").Value = 0.1192928843907551
$ws.Range("E6").Value = 0.3453880200452168
$ws.Range("F6").Value = 0.3463346331599631
$ws.Range("G6").Value = 19

$ws.Range("B7").Value = 0.02628598706572433
$ws.Range("C7").Value = 0.2431904892130187
$ws.Range("D7").Value = 0.0933406381209692
$ws.Range("E7").Value = 0.3055170013615759
$ws.Range("F7").Value = 0.3132086628365013
$ws.Range("G7").Value = 18

$ws.Range("B8").Value = 0.05644726556593638
$ws.Range("C8").Value = 0.2954772732059698
$ws.Range("D8").Value = 0.1339965475761727
$ws.Range("E8").Value = 0.3660553886724968
$ws.Range("F8").Value = 0.372808120415778
$ws.Range("G8").Value = 17
